# Update the "Project list" sheet: add a "Domain" column (inserted between
# "Project name" and "Description") and fill in the first data row with the
# Credit Risk Modeling project.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlCenter = -4108   # xlCenter / xlHAlignCenter / xlVAlignCenter
$xlTop    = -4160   # xlTop / xlVAlignTop

# Insert a new column for "Domain" right after "Project name" (old column C
# "Description" and everything after it shifts one column right).
$ws.Columns("C").Insert()
$ws.Columns("C").ColumnWidth = $ws.Columns("B").ColumnWidth

# --- Header row ---
$ws.Range("C1").Value = "Domain"

# --- First data row ---
$ws.Range("A2").Value = 1
$ws.Range("A2").HorizontalAlignment = $xlCenter
$ws.Range("A2").VerticalAlignment = $xlTop

$ws.Range("C2").Value = "Finance"
$ws.Range("C2").VerticalAlignment = $xlTop

$ws.Range("B2").Value = "Credit Risk Modeling"
$ws.Range("B2").VerticalAlignment = $xlTop

$ws.Range("D2").Value = "Model to predict if the loan provided by a financial institution will default or not. The output generated is a probability of failing to repay the loan. Higher the value, greater the risk of defaulting."
$ws.Range("D2").VerticalAlignment = $xlTop
$ws.Range("D2").WrapText = $true

$ws.Range("F2").Value = "In-progress"
$ws.Range("F2").VerticalAlignment = $xlTop

# Row grows tall to fit the wrapped description text.
$ws.Rows("2").RowHeight = 57.6

# Leave selection where the author left it.
$ws.Range("E2").Select()
